# Update the 25 multiplication-problem answers shown in the table cells
# of the generated worksheet to match the newly regenerated values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("37×49=1813", $true, $false, $false, $false, $false, $true, 1, $false, "27×79=2133", 2) | Out-Null
$d.Content.Find.Execute("69×49=3381", $true, $false, $false, $false, $false, $true, 1, $false, "77×53=4081", 2) | Out-Null
$d.Content.Find.Execute("49×82=4018", $true, $false, $false, $false, $false, $true, 1, $false, "61×93=5673", 2) | Out-Null
$d.Content.Find.Execute("45×21=945", $true, $false, $false, $false, $false, $true, 1, $false, "62×66=4092", 2) | Out-Null
$d.Content.Find.Execute("95×36=3420", $true, $false, $false, $false, $false, $true, 1, $false, "23×65=1495", 2) | Out-Null
$d.Content.Find.Execute("99×40=3960", $true, $false, $false, $false, $false, $true, 1, $false, "26×33=858", 2) | Out-Null
$d.Content.Find.Execute("40×95=3800", $true, $false, $false, $false, $false, $true, 1, $false, "63×60=3780", 2) | Out-Null
$d.Content.Find.Execute("79×54=4266", $true, $false, $false, $false, $false, $true, 1, $false, "15×41=615", 2) | Out-Null
$d.Content.Find.Execute("98×71=6958", $true, $false, $false, $false, $false, $true, 1, $false, "71×18=1278", 2) | Out-Null
$d.Content.Find.Execute("50×85=4250", $true, $false, $false, $false, $false, $true, 1, $false, "13×74=962", 2) | Out-Null
$d.Content.Find.Execute("28×25=700", $true, $false, $false, $false, $false, $true, 1, $false, "69×83=5727", 2) | Out-Null
$d.Content.Find.Execute("37×28=1036", $true, $false, $false, $false, $false, $true, 1, $false, "92×34=3128", 2) | Out-Null
$d.Content.Find.Execute("34×38=1292", $true, $false, $false, $false, $false, $true, 1, $false, "89×41=3649", 2) | Out-Null
$d.Content.Find.Execute("75×79=5925", $true, $false, $false, $false, $false, $true, 1, $false, "98×48=4704", 2) | Out-Null
$d.Content.Find.Execute("96×42=4032", $true, $false, $false, $false, $false, $true, 1, $false, "42×67=2814", 2) | Out-Null
$d.Content.Find.Execute("60×13=780", $true, $false, $false, $false, $false, $true, 1, $false, "12×60=720", 2) | Out-Null
$d.Content.Find.Execute("43×44=1892", $true, $false, $false, $false, $false, $true, 1, $false, "83×97=8051", 2) | Out-Null
$d.Content.Find.Execute("53×78=4134", $true, $false, $false, $false, $false, $true, 1, $false, "87×47=4089", 2) | Out-Null
$d.Content.Find.Execute("46×71=3266", $true, $false, $false, $false, $false, $true, 1, $false, "28×80=2240", 2) | Out-Null
$d.Content.Find.Execute("92×67=6164", $true, $false, $false, $false, $false, $true, 1, $false, "77×22=1694", 2) | Out-Null
$d.Content.Find.Execute("86×24=2064", $true, $false, $false, $false, $false, $true, 1, $false, "76×66=5016", 2) | Out-Null
$d.Content.Find.Execute("27×30=810", $true, $false, $false, $false, $false, $true, 1, $false, "24×30=720", 2) | Out-Null
$d.Content.Find.Execute("45×82=3690", $true, $false, $false, $false, $false, $true, 1, $false, "72×41=2952", 2) | Out-Null
$d.Content.Find.Execute("59×66=3894", $true, $false, $false, $false, $false, $true, 1, $false, "67×20=1340", 2) | Out-Null
$d.Content.Find.Execute("81×76=6156", $true, $false, $false, $false, $false, $true, 1, $false, "13×61=793", 2) | Out-Null
